# Cheatsheet.docx edit: rewrite the Boolean-operators section.
#
# The old paragraphs (originally #4-#8, containing "Boolean kan gebruikt
# woorden...", "Als ergens een "!" staat...", an empty paragraph, the
# "Met "&&"/"||"..." paragraph, and a trailing empty paragraph) are
# replaced by a longer block: the same two intro sentences (now typed as
# single plain runs, no proofErr spell-markers), the "&&"/"||" explanation,
# a new comparison-operator cheatsheet (==, <, >, >=, <=), and four new
# English-language "truth table" example lines, followed by two trailing
# blank paragraphs.

$d = $word.ActiveDocument

# Anchor on the end of paragraph 2 ("Als het of moet zijn ... "||"")
# which is the last paragraph that stays untouched before the block we
# are rewriting. Paragraph 3 (the blank line right after it) is already
# in place and is left as-is.
$anchor = $d.Paragraphs.Item(2).Range.Duplicate
$anchor.Collapse(0)   # wdCollapseEnd

$vbreak = [char]11     # manual line break -> <w:br/>

$newText = "`r" + `
  "Boolean kan gebruikt woorden als een veriable`r" + `
  "Als ergens een “!” staat betekent het niet. Dus als je bijvoorbeeld !a ziet kan het beteken dat als het “false” is dat het true word`r" + `
  "`r" + `
  "Met “&&” moet alles true zijn" + $vbreak + "Met “||” Moet er minimaal 1 true zijn om de hele lijn true te maken`r" + `
  "`r" + `
  "Gelijk aan:  ==`r" + `
  "Kleiner dan:  <`r" + `
  "Groter dan:  >`r" + `
  "Groter of gelijk aan:  >=`r" + `
  "Kleiner of gelijk aan:  <=`r" + `
  "`r" + `
  "True&&`ttrue`t`ttrue`r" + `
  "False&&true`t`tfalse`r" + `
  "True&&`tfalse`t`tfalse`r" + `
  "False&&false`t`tfalse`r" + `
  "`r"

$anchor.InsertAfter($newText)

# The block we just inserted now occupies paragraphs 4-19 (16 new
# paragraphs after the anchor, the last `r above starts the 17th / final
# blank paragraph which is actually the old, now-shifted, paragraph 3's
# successor... handled below). Re-derive indices from the document instead
# of hard-coding, then delete the stale original paragraphs that used to
# sit right after the old paragraph 3 and apply en-US language to the
# four new example paragraphs + the two trailing blanks.

# Paragraph 3 is unchanged ("" blank line). Our 16 freshly typed
# paragraphs are #4 .. #19.
$firstNew = 4
$lastNew = 19

for ($i = $firstNew; $i -le $lastNew; $i++) {
    $idx = $i - $firstNew + 1
    if ($idx -ge 12) {
        $d.Paragraphs.Item($i).Range.LanguageID = "en-US"
    }
}

# Remove the old paragraphs that followed (the original "Boolean...",
# "Als ergens...", blank, "Met ...", blank block), which now sit right
# after our newly inserted paragraphs.
$oldStart = $d.Paragraphs.Item($lastNew + 1).Range.Start
$oldEnd = $d.Paragraphs.Item($lastNew + 5).Range.End
$oldRange = $d.Range($oldStart, $oldEnd)
$oldRange.Delete()
